# Add 2022-Q3 data
# 1. Update the "总计" (summary) sheet: insert a new leading row for 2022-Q3
#    and shift the existing rows down, adding the 2020-Q4 row at the end.
# 2. Insert a brand-new "2022-Q3" worksheet (holding per-fund data) right
#    before the existing "2022-Q2" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the new last row (row 7) the same number formatting as row 6 before
# overwriting its values (keeps the bold-index style on column A).
$total.Cells.Item(6,1).Copy($total.Cells.Item(7,1))

# index, quarter label, 持有数量(只), 持有市值(亿元) - rows 2..7 after the edit
$totalRows = @(
    @(0, "2022-Q3", 5, 5.44),
    @(1, "2022-Q2", 3, 2.74),
    @(2, "2022-Q1", 9, 4.73),
    @(3, "2021-Q4", 2, 5.93),
    @(4, "2021-Q3", 7, 5.07),
    @(5, "2020-Q4", 1, 0.05)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r,1).Value2 = $row[0]
    $total.Cells.Item($r,2).Value2 = $row[1]
    $total.Cells.Item($r,3).Value2 = $row[2]
    $total.Cells.Item($r,4).Value2 = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with per-fund holdings
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Pull formatting (bold header / bold index column) from the "总计" sheet so
# the new tab matches the look of its siblings without dragging any values.
$total.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2:A6").Copy()
$q3.Range("A2:A6").PasteSpecial(-4122)

$header = New-Object 'object[,]' 1,7
$header[0,0] = "基金代码"
$header[0,1] = "基金名称"
$header[0,2] = "基金规模"
$header[0,3] = "股票总仓位"
$header[0,4] = "仓位占比"
$header[0,5] = "持有市值(亿元)"
$header[0,6] = "仓位排名"
$q3.Range("B1:H1").Value2 = $header

# index, fund code, fund name, scale, total stock position, position ratio,
# market value held (亿元), position rank
$rows = @(
    @(0, "'010363", "信澳匠心臻选两年持有期混合",   "'48.09", "'92.07", "'4.92", "'2.3660", 3),
    @(1, "'010963", "信澳周期动力混合A",            "'39.77", "'89.47", "'5.91", "'2.3504", 1),
    @(2, "'015455", "信澳周期动力混合C",            "'12.28", "'89.47", "'5.91", "'0.7257", 1),
    @(3, "'003587", "先锋精一灵活配置混合C",         "'0.02",  "'94.32", "'2.64", "'0.0005", 5),
    @(4, "'003586", "先锋精一灵活配置混合A",         "'0.01",  "'94.32", "'2.64", "'0.0003", 5)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r,1).Value2 = $row[0]
    $data = New-Object 'object[,]' 1,7
    for ($i = 0; $i -lt 7; $i++) { $data[0,$i] = $row[$i + 1] }
    $q3.Range("B$r`:H$r").Value2 = $data
    $r = $r + 1
}

# The leading apostrophes above force Excel to keep numeric-looking text
# (fund codes, percentages, ...) as text, same as the source data - but they
# also stamp a "quote prefix" style on those cells. Strip that back off so
# the cells end up with plain (unstyled) text, matching the rest of the
# workbook.
$q3.Range("B2:G6").ClearFormats()

Write-Host "2022-Q3 sheet added and 总计 sheet updated"
